$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.235.99"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").Value = "3.489.94"
$ws.Range("E3").Value = "  +2.27%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +9.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("E12").Value = "  +4.65%  "

$ws.Range("E13").Value = "  +4.37%  "

$ws.Range("D14").Value = "4.045.01"
$ws.Range("E14").Value = "  +2.42%  "

$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.47%  "

$ws.Range("D17").Value = "3.480.14"
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").Value = "63.082.75"
$ws.Range("E20").Value = "  +1.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  +3.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.50%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.43%  "

$ws.Range("E37").Value = "  -2.97%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.80%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.85%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.135"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.90%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "148.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.01%  "

$ws.Range("E43").Value = "  +1.66%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.70%  "

$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("E46").Value = "  +2.73%  "

$ws.Range("D47").Value = "0.0₃0592"
$ws.Range("E47").Value = "  +38.08%  "

$ws.Range("E48").Value = "  +9.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.146"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "

